$wb = $excel.ActiveWorkbook
$wsLugar = $wb.Worksheets.Item("Lugar")
$wsHist  = $wb.Worksheets.Item("Hitorial")

# --- Hitorial: insert a new "Estado" column before the existing "Imagen"
#     column (old E -> new F), then fill in the header + O/C status values.
$wsHist.Columns.Item(5).Insert()

$wsHist.Range("E1").Value = "Estado"
$wsHist.Range("E2").Value = "O"
$wsHist.Range("E3").Value = "C"
$wsHist.Range("E4").Value = "O"
$wsHist.Range("E5").Value = "O"

# The column insert shifted the Imagen values/hyperlink styling to column F
# automatically, but the <hyperlinks> entries themselves still reference the
# old column E addresses, so rebuild them against column F.
$wsHist.Hyperlinks.Delete()
$wsHist.Hyperlinks.Add($wsHist.Range("F3"), "https://raw.githubusercontent.com/brauliovargas/Huisachito/master/img/Venado02.jpg")
$wsHist.Hyperlinks.Add($wsHist.Range("F4"), "https://raw.githubusercontent.com/brauliovargas/Huisachito/master/img/Venado03.jpg")
$wsHist.Hyperlinks.Add($wsHist.Range("F5"), "https://raw.githubusercontent.com/brauliovargas/Huisachito/master/img/Venado04.jpg")

# Re-apply the workbook's hyperlink cell style (Add() resets it) so the
# linked cells keep looking like the rest of the hyperlinked column.
$wsHist.Range("F3:F5").Style = "Hipervínculo"

# --- Tab / selection state: "Hitorial" becomes the active sheet (was
#     "Lugar"), with the cursor left on E6; "Lugar" keeps its A9 selection.
$wsLugar.Range("A9").Select()
$wsHist.Range("E6").Select()

Write-Output "ok"
